# Apply trade #42 close update across the workbook.
$wb = $excel.ActiveWorkbook

# --- Summary sheet updates ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.8   # Current Capital
$summary.Range("B4").Value = 0.8      # Total P&L $
$summary.Range("B5").Value = 0.38     # Total P&L %
$summary.Range("B6").Value = 42       # Total Trades
$summary.Range("B7").Value = 14       # Winning Trades
$summary.Range("B9").Value = 33.33    # Win Rate %

# --- Strategy Status sheet updates (MarketMaking row, row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.8     # Capital
$status.Range("D4").Value = 42        # Trades
$status.Range("E4").Value = 0.8       # P&L $
$status.Range("F4").Value = 0.8       # P&L %
$status.Range("G4").Value = 33.33     # Win Rate %

# --- New trade row (row 43) appended to "All Trades" and "MarketMaking" sheets ---
$newRow = @{
    A = 42
    B = "2026-02-17"
    C = "15:29:00"
    D = "MarketMaking"
    E = "UP"
    F = 0.03
    G = 0.58
    H = "CLOSED"
    I = 1833.3333
    J = 0.55
    K = 100.8
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 1.14
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A43").Value = $newRow.A
    # Date/Time text columns: force text so Excel doesn't auto-convert the
    # string into a date/time serial number (matches the other rows, which
    # store the date/time as plain text, not as real dates).
    $ws.Range("B43").NumberFormat = "@"
    $ws.Range("B43").Value = $newRow.B
    $ws.Range("C43").NumberFormat = "@"
    $ws.Range("C43").Value = $newRow.C
    $ws.Range("D43").Value = $newRow.D
    $ws.Range("E43").Value = $newRow.E
    $ws.Range("F43").Value = $newRow.F
    $ws.Range("G43").Value = $newRow.G
    $ws.Range("H43").Value = $newRow.H
    $ws.Range("I43").Value = $newRow.I
    $ws.Range("J43").Value = $newRow.J
    $ws.Range("K43").Value = $newRow.K
    $ws.Range("L43").Value = $newRow.L
    $ws.Range("M43").Value = $newRow.M
    $ws.Range("N43").Value = $newRow.N
    $ws.Range("O43").Value = $newRow.O
    $ws.Range("P43").Value = $newRow.P
    $ws.Range("Q43").Value = $newRow.Q
}

Write-Host "Done applying trade #42 close update."
